# Update workbook version string from
#   "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
# to
#   "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"
# across the "About" sheet and the "Boundaries and methane sources" sheet.

$wb = $excel.ActiveWorkbook

$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- "About" sheet ---
$aboutSheet = $wb.Worksheets.Item("About")

$aboutSheet.Range("A2").Value = "Version: " + $newVersion

$aboutSheet.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Shoushan No.1 Coal Mine, China, M2131, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 8; $row++) {
    $dataSheet.Range("S" + $row).Value = $newVersion
}
